$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values are plain decimal-looking strings (e.g. "516.47") that
# Excel would otherwise auto-convert to numbers, losing exact text formatting
# and introducing floating point noise. Force those cells to Text format first
# so the assigned values round-trip as the exact original strings.
$textForceCells = @(
    "D5",
    "D6",
    "D7",
    "D10",
    "D11",
    "D13",
    "D14",
    "D16",
    "D18",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D38",
    "D39",
    "D40",
    "D42",
    "D44",
    "D45",
    "D47",
    "D49",
    "D50",
    "D51"
)
foreach ($ref in $textForceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '61.898.92'
$ws.Range("E2").Value = '  -4.75%  '

# Row 3
$ws.Range("D3").Value = '3.147.05'
$ws.Range("E3").Value = '  -6.64%  '

# Row 4
$ws.Range("E4").Value = '  +0.17%  '

# Row 5
$ws.Range("D5").Value = '516.47'
$ws.Range("E5").Value = '  -2.19%  '

# Row 6
$ws.Range("D6").Value = '168.20'
$ws.Range("E6").Value = '  -8.40%  '

# Row 7
$ws.Range("D7").Value = '0.584'
$ws.Range("E7").Value = '  -2.18%  '

# Row 8
$ws.Range("E8").Value = '  +0.04%  '

# Row 9
$ws.Range("D9").Value = '3.147.58'
$ws.Range("E9").Value = '  -6.51%  '

# Row 10
$ws.Range("D10").Value = '0.592'
$ws.Range("E10").Value = '  -4.20%  '

# Row 11
$ws.Range("D11").Value = '51.70'
$ws.Range("E11").Value = '  -9.32%  '

# Row 12
$ws.Range("E12").Value = '  -2.76%  '

# Row 13
$ws.Range("D13").Value = '0.0000245'
$ws.Range("E13").Value = '  -2.72%  '

# Row 14
$ws.Range("D14").Value = '8.85'
$ws.Range("E14").Value = '  -3.27%  '

# Row 15
$ws.Range("D15").Value = '3.640.61'
$ws.Range("E15").Value = '  -6.28%  '

# Row 16
$ws.Range("D16").Value = '0.115'
$ws.Range("E16").Value = '  -5.31%  '

# Row 17
$ws.Range("D17").Value = '3.140.56'
$ws.Range("E17").Value = '  -6.46%  '

# Row 18
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").Value = '16.89'
$ws.Range("E18").Value = '  -2.26%  '

# Row 19
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '61.731.51'
$ws.Range("E19").Value = '  -4.39%  '

# Row 20
$ws.Range("D20").Value = '0.958'
$ws.Range("E20").Value = '  -0.09%  '

# Row 21
$ws.Range("D21").Value = '10.78'
$ws.Range("E21").Value = '  -1.62%  '

# Row 22
$ws.Range("D22").Value = '357.99'
$ws.Range("E22").Value = '  -3.53%  '

# Row 23
$ws.Range("D23").Value = '11.06'
$ws.Range("E23").Value = '  +3.47%  '

# Row 24
$ws.Range("D24").Value = '3.66'
$ws.Range("E24").Value = '  -0.99%  '

# Row 25
$ws.Range("D25").Value = '79.45'
$ws.Range("E25").Value = '  -1.56%  '

# Row 26
$ws.Range("D26").Value = '3.88'
$ws.Range("E26").Value = '  +3.80%  '

# Row 27
$ws.Range("D27").Value = '6.11'
$ws.Range("E27").Value = '  +4.11%  '

# Row 28
$ws.Range("D28").Value = '2.57'
$ws.Range("E28").Value = '  -2.02%  '

# Row 29
$ws.Range("D29").Value = '11.06'
$ws.Range("E29").Value = '  -1.57%  '

# Row 30
$ws.Range("D30").Value = '7.99'
$ws.Range("E30").Value = '  -4.35%  '

# Row 31
$ws.Range("D31").Value = '632.47'
$ws.Range("E31").Value = '  -5.16%  '

# Row 32
$ws.Range("D32").Value = '27.72'
$ws.Range("E32").Value = '  -4.51%  '

# Row 33
$ws.Range("D33").Value = '6.31'
$ws.Range("E33").Value = '  -6.11%  '

# Row 34
$ws.Range("D34").Value = '11.14'
$ws.Range("E34").Value = '  +0.85%  '

# Row 35
$ws.Range("D35").Value = '0.103'
$ws.Range("E35").Value = '  -0.79%  '

# Row 36
$ws.Range("D36").Value = '55.99'
$ws.Range("E36").Value = '  -6.89%  '

# Row 37
$ws.Range("E37").Value = '  -0.01%  '

# Row 38
$ws.Range("D38").Value = '36.32'
$ws.Range("E38").Value = '  +0.54%  '

# Row 39
$ws.Range("D39").Value = '0.367'
$ws.Range("E39").Value = '  -1.87%  '

# Row 40
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  +0.26%  '

# Row 41
$ws.Range("D41").Value = '0.0₃0684'
$ws.Range("E41").Value = '  +11.68%  '

# Row 42
$ws.Range("D42").Value = '0.123'
$ws.Range("E42").Value = '  -2.78%  '

# Row 43
$ws.Range("D43").Value = '2.874.28'
$ws.Range("E43").Value = '  +2.88%  '

# Row 44
$ws.Range("D44").Value = '2.49'
$ws.Range("E44").Value = '  +8.57%  '

# Row 45
$ws.Range("D45").Value = '2.87'
$ws.Range("E45").Value = '  +11.96%  '

# Row 46
$ws.Range("E46").Value = '  +0.77%  '

# Row 47
$ws.Range("D47").Value = '0.0383'
$ws.Range("E47").Value = '  -0.30%  '

# Row 48
$ws.Range("E48").Value = '  +4.65%  '

# Row 49
$ws.Range("D49").Value = '2.51'
$ws.Range("E49").Value = '  -8.39%  '

# Row 50
$ws.Range("D50").Value = '0.122'
$ws.Range("E50").Value = '  -2.41%  '

# Row 51
$ws.Range("D51").Value = '132.18'
$ws.Range("E51").Value = '  -2.50%  '
